$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Description" note for the pRSFduet primers (rows 4 and 5).
$ws.Range("E4").Value = "pRSF vector 증폭용 primer"
$ws.Range("E5").Value = "pRSF vector 증폭용 primer"

# Remove the now-unused scratch rows (old position placeholders A9:A12, B1:B12)
# that used to live below the primer table.
$ws.Rows("10:25").Delete()

# Widen column E so the new description text is readable.
$ws.Columns("E").ColumnWidth = 39.428571428571427

# Match the author's last on-screen selection.
$ws.Range("C16").Select()
